# The deck's theme was swapped: the "Default" color palette (previously
# only used by the notes master's theme part) becomes the presentation's
# main theme (the slide master's theme), while the "Simple Light" palette
# (previously the main theme) moves off the slide master.
#
# PowerPoint's ColorFormat.RGB is a VBA-style BGR-packed long (0xBBGGRR),
# so convert each target RRGGBB hex value from the OOXML <a:srgbClr val=.../>
# before assigning it to ThemeColorScheme.Colors(i).RGB.

function Convert-HexToBgrLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target ("Default") theme color scheme, in clrScheme schema order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$newThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "158158",  # dk2
    "F3F3F3",  # lt2
    "058DC7",  # accent1
    "50B432",  # accent2
    "ED561B",  # accent3
    "EDEF00",  # accent4
    "24CBE5",  # accent5
    "64E572",  # accent6
    "2200CC",  # hlink
    "551A8B"   # folHlink
)

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$themeColorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $newThemeColors.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = Convert-HexToBgrLong $newThemeColors[$i - 1]
}
